# This script applies the IFRS data-correction edit described in the commit
# "error solve ifrs list": rows 2-6 get corrected numeric values (D:AJ),
# a few now-blank cells (AG/AH on row 5, AG/AH/AI on row 6) are cleared,
# and rows 7-9 have all of their data columns (D:AI) cleared, leaving only
# the leading A/B/C identifier columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5069
$ws.Range("E2").Value = -214
$ws.Range("F2").Value = -160
$ws.Range("G2").Value = -843
$ws.Range("H2").Value = -744
$ws.Range("I2").Value = -732
$ws.Range("J2").Value = -11
$ws.Range("K2").Value = 7262
$ws.Range("L2").Value = 6445
$ws.Range("M2").Value = 817
$ws.Range("N2").Value = 697
$ws.Range("O2").Value = 121
$ws.Range("P2").Value = 326
$ws.Range("Q2").Value = -24
$ws.Range("R2").Value = -306
$ws.Range("S2").Value = 332
$ws.Range("T2").Value = 237
$ws.Range("U2").Value = -261
$ws.Range("V2").Value = 5443
$ws.Range("W2").Value = -4.21
$ws.Range("X2").Value = -14.67
$ws.Range("Y2").Value = -69.15000000000001
$ws.Range("Z2").Value = -10.26
$ws.Range("AA2").Value = 788.52
$ws.Range("AB2").Value = 211.47
$ws.Range("AC2").Value = -1124
$ws.Range("AD2").Value = -2.63
$ws.Range("AE2").Value = 1144
$ws.Range("AF2").Value = 2.58
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 65180349

# Row 3
$ws.Range("D3").Value = 4914
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 18
$ws.Range("G3").Value = -406
$ws.Range("H3").Value = -346
$ws.Range("I3").Value = -331
$ws.Range("J3").Value = -15
$ws.Range("K3").Value = 5393
$ws.Range("L3").Value = 4780
$ws.Range("M3").Value = 613
$ws.Range("N3").Value = 497
$ws.Range("O3").Value = 116
$ws.Range("P3").Value = 326
$ws.Range("Q3").Value = 132
$ws.Range("R3").Value = 548
$ws.Range("S3").Value = -685
$ws.Range("T3").Value = 135
$ws.Range("U3").Value = -4
$ws.Range("V3").Value = 4031
$ws.Range("W3").Value = 0.63
$ws.Range("X3").Value = -7.04
$ws.Range("Y3").Value = -55.46
$ws.Range("Z3").Value = -5.47
$ws.Range("AA3").Value = 779.46
$ws.Range("AB3").Value = 106.95
$ws.Range("AC3").Value = -508
$ws.Range("AD3").Value = -3.39
$ws.Range("AE3").Value = 775
$ws.Range("AF3").Value = 2.22
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 65180349

# Row 4
$ws.Range("D4").Value = 3818
$ws.Range("E4").Value = 259
$ws.Range("F4").Value = 265
$ws.Range("G4").Value = 69
$ws.Range("H4").Value = 48
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 4675
$ws.Range("L4").Value = 3956
$ws.Range("M4").Value = 719
$ws.Range("N4").Value = 592
$ws.Range("O4").Value = 127
$ws.Range("P4").Value = 326
$ws.Range("Q4").Value = 469
$ws.Range("R4").Value = 293
$ws.Range("S4").Value = -769
$ws.Range("T4").Value = 61
$ws.Range("U4").Value = 409
$ws.Range("V4").Value = 3262
$ws.Range("W4").Value = 6.77
$ws.Range("X4").Value = 1.27
$ws.Range("Y4").Value = 7.98
$ws.Range("Z4").Value = 0.96
$ws.Range("AA4").Value = 550.3099999999999
$ws.Range("AB4").Value = 122.44
$ws.Range("AC4").Value = 67
$ws.Range("AD4").Value = 24.75
$ws.Range("AE4").Value = 909
$ws.Range("AF4").Value = 1.82
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 65180349

# Row 5 (AG5 and AH5 become blank - no value at all)
$ws.Range("D5").Value = 3951
$ws.Range("E5").Value = 293
$ws.Range("F5").Value = 293
$ws.Range("G5").Value = 360
$ws.Range("H5").Value = 292
$ws.Range("I5").Value = 293
$ws.Range("J5").Value = -1
$ws.Range("K5").Value = 4898
$ws.Range("L5").Value = 3278
$ws.Range("M5").Value = 1621
$ws.Range("N5").Value = 1627
$ws.Range("O5").Value = -6
$ws.Range("P5").Value = 706
$ws.Range("Q5").Value = 298
$ws.Range("R5").Value = 283
$ws.Range("S5").Value = -567
$ws.Range("T5").Value = 32
$ws.Range("U5").Value = 266
$ws.Range("V5").Value = 2731
$ws.Range("W5").Value = 7.41
$ws.Range("X5").Value = 7.39
$ws.Range("Y5").Value = 26.41
$ws.Range("Z5").Value = 6.1
$ws.Range("AA5").Value = 202.25
$ws.Range("AB5").Value = 175.53
$ws.Range("AC5").Value = 323
$ws.Range("AD5").Value = 4.37
$ws.Range("AE5").Value = 1379
$ws.Range("AF5").Value = 1.02
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 141144600

# Row 6 (AG6, AH6, AI6 become blank - no value at all; J6/O6 already blank)
$ws.Range("D6").Value = 4225
$ws.Range("E6").Value = 258
$ws.Range("F6").Value = 258
$ws.Range("G6").Value = 155
$ws.Range("H6").Value = 157
$ws.Range("I6").Value = 159
$ws.Range("K6").Value = 4682
$ws.Range("L6").Value = 2910
$ws.Range("M6").Value = 1772
$ws.Range("N6").Value = 1781
$ws.Range("P6").Value = 706
$ws.Range("Q6").Value = 287
$ws.Range("R6").Value = -18
$ws.Range("S6").Value = -298
$ws.Range("T6").Value = 38
$ws.Range("U6").Value = 249
$ws.Range("V6").Value = 2471
$ws.Range("W6").Value = 6.11
$ws.Range("X6").Value = 3.72
$ws.Range("Y6").Value = 9.32
$ws.Range("Z6").Value = 3.28
$ws.Range("AA6").Value = 164.29
$ws.Range("AB6").Value = 200.14
$ws.Range("AC6").Value = 113
$ws.Range("AD6").Value = 10.53
$ws.Range("AE6").Value = 1509
$ws.Range("AF6").Value = 0.79
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()
$ws.Range("AJ6").Value = 141144600

# Rows 7-9: entire data portion (D:AI) is cleared, leaving only A/B/C populated
$ws.Range("D7:AI9").ClearContents()

Write-Output "edit applied"
